$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The Python/IDLE version in the Software Versions table was listed
# incorrectly; fix it from 2.8.3rc1 to 3.8.3rc1.
$ws.Range("B4").Value = "3.8.3rc1"
